# feat: add 2022-Q3 data
#
# 1) "总计" (summary) sheet: the old row for 2022-Q2 becomes the new
#    2022-Q3 totals, a fresh row is inserted for (the now-historical)
#    2022-Q2 totals, and the 2022-Q1 row shifts down one row.
# 2) A brand-new "2022-Q3" worksheet is inserted right after "总计" and
#    before "2022-Q2", holding the per-fund holdings detail for the
#    new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: update the "总计" (totals) sheet
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Row 2 used to be "2022-Q2" -> becomes the new "2022-Q3" totals
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.08

# Duplicate the formatting of the A-column index cell down into the two
# rows below it, so the new/shifted rows keep the same style (s=2) as
# row 2 without introducing any new style entries.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3:A4").PasteSpecial(-4122)

# New row 3: the historical "2022-Q2" totals (what used to live in row 2)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.02

# Row 4: "2022-Q1" totals, shifted down from the old row 3
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 4
$wsTotal.Range("D4").Value = 0.35

# ---------------------------------------------------------------------
# Part 2: insert the new "2022-Q3" detail sheet
# ---------------------------------------------------------------------
# Duplicate the existing "2022-Q2" sheet (position 2) so the new sheet
# starts out with identical column layout/styling, and place the copy
# immediately before it -> it lands in position 2, "2022-Q2" shifts to
# position 3, "2022-Q1" shifts to position 4.
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($wsQ2)

$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# The fund-code column (B) and the percentage-ish numeric-looking text
# columns (D:G) must stay plain text (e.g. "006165" keeps its leading
# zero, "3.87" stays a string rather than becoming the number 3.87) -
# mark them as Text before writing the values.
$wsQ3.Range("B2:B4").NumberFormat = "@"
$wsQ3.Range("D2:G4").NumberFormat = "@"

# Row 2 - 建信中证1000指数增强A
$wsQ3.Range("B2").Value = "006165"
$wsQ3.Range("C2").Value = "建信中证1000指数增强A"
$wsQ3.Range("D2").Value = "3.87"
$wsQ3.Range("E2").Value = "84.02"
$wsQ3.Range("F2").Value = "1.35"
$wsQ3.Range("G2").Value = "0.0522"
$wsQ3.Range("H2").Value = 6

# Row 3 - 建信中证1000指数增强C
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "006166"
$wsQ3.Range("C3").Value = "建信中证1000指数增强C"
$wsQ3.Range("D3").Value = "1.89"
$wsQ3.Range("E3").Value = "84.02"
$wsQ3.Range("F3").Value = "1.35"
$wsQ3.Range("G3").Value = "0.0255"
$wsQ3.Range("H3").Value = 6

# Row 4 - 建信中证1000指数增强E
$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").Value = "013442"
$wsQ3.Range("C4").Value = "建信中证1000指数增强E"
$wsQ3.Range("D4").Value = "0.18"
$wsQ3.Range("E4").Value = "84.02"
$wsQ3.Range("F4").Value = "1.35"
$wsQ3.Range("G4").Value = "0.0024"
$wsQ3.Range("H4").Value = 6

# Give row 3/4's index column (A) the same style as row 2's (copied from
# the duplicated "2022-Q2" sheet), same trick as on the totals sheet.
$wsQ3.Range("A2").Copy()
$wsQ3.Range("A3:A4").PasteSpecial(-4122)
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("A4").Value = 2

# Leave the workbook showing the "总计" sheet, same as before the edit.
$wsTotal.Activate()
